$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45192 -> 2023-09-23).
# Update every data row (2 through 111) to the new date serial 45202 (2023-10-03).
$oldSerial = 45192
$newSerial = 45202

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 111 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
